$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
# si#6 "Volume 32   Number  14" -> "...15"; si#9 date range 3/31-4/6 -> 4/7-4/13
$ws.Range("A8").Value = "Volume 32   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/7/2025  Through  4/13/2025"

# --- Cells whose type/style switches between numeric and text placeholder ---
# Copy a same-row/col donor cell first so the destination picks up the correct
# style + shared-text value, then overwrite with the real value where needed.
$ws.Range("L15").Copy($ws.Range("L14"))       # "***.*" text -> numeric -100 (style 15)
$ws.Range("L14").Value = -100

$ws.Range("D22").Copy($ws.Range("C22"))       # numeric 1 -> text "0" (style 13)
$ws.Range("D27").Copy($ws.Range("C27"))       # numeric 1 -> text "0" (style 13)

$ws.Range("F29").Copy($ws.Range("D29"))       # text "0" -> numeric (style 14)
$ws.Range("D29").Value = 1
$ws.Range("K29").Copy($ws.Range("E29"))       # text "***.*" -> numeric (style 15)
$ws.Range("E29").Value = -100

$ws.Range("F30").Copy($ws.Range("D30"))       # text "0" -> numeric (style 14)
$ws.Range("D30").Value = 1
$ws.Range("K30").Copy($ws.Range("E30"))       # text "***.*" -> numeric (style 15)
$ws.Range("E30").Value = -100

# --- Updated weekly crime statistics values ---
# Row 15
$ws.Range("F15").Value = 1
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = -14.285714285714
# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("F16").Value = 13
$ws.Range("H16").Value = -13.333333333333
$ws.Range("I16").Value = 53
$ws.Range("J16").Value = 68
$ws.Range("K16").Value = -22.058823529411
$ws.Range("L16").Value = -23.188405797101
$ws.Range("M16").Value = -17.1875
$ws.Range("N16").Value = -73.631840796019
# Row 17
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -36.666666666666
$ws.Range("I17").Value = 72
$ws.Range("J17").Value = 125
$ws.Range("K17").Value = -42.4
$ws.Range("L17").Value = -20.879120879120
$ws.Range("M17").Value = 12.5
$ws.Range("N17").Value = -54.716981132075
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 12.5
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 28
$ws.Range("L18").Value = 3.225806451612
$ws.Range("M18").Value = 6.666666666666
$ws.Range("N18").Value = -75.757575757575
# Row 19
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -50
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -29.411764705882
$ws.Range("I19").Value = 105
$ws.Range("J19").Value = 139
$ws.Range("K19").Value = -24.460431654676
$ws.Range("L19").Value = -3.669724770642
$ws.Range("M19").Value = 94.444444444444
$ws.Range("N19").Value = 34.615384615384
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -57.142857142857
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = -51.428571428571
$ws.Range("L20").Value = 21.428571428571
$ws.Range("M20").Value = 41.666666666666
$ws.Range("N20").Value = -78.75
# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -37.5
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 122
$ws.Range("H21").Value = -31.147540983606
$ws.Range("I21").Value = 285
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = -28.75
$ws.Range("L21").Value = -10.9375
$ws.Range("M21").Value = 21.794871794871
$ws.Range("N21").Value = -56.818181818181
# Row 22
$ws.Range("L22").Value = -37.5
$ws.Range("M22").Value = -68.75
# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -9.090909090909
$ws.Range("I23").Value = 42
$ws.Range("J23").Value = 53
$ws.Range("K23").Value = -20.754716981132
$ws.Range("L23").Value = -4.545454545454
$ws.Range("M23").Value = 90.909090909090
# Row 24
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -47.619047619047
$ws.Range("F24").Value = 49
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -40.963855421686
$ws.Range("I24").Value = 213
$ws.Range("J24").Value = 256
$ws.Range("K24").Value = -16.796875
$ws.Range("L24").Value = -40.833333333333
$ws.Range("M24").Value = -20.522388059701
# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -20
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = -44
$ws.Range("I25").Value = 59
$ws.Range("J25").Value = 68
$ws.Range("K25").Value = -13.235294117647
$ws.Range("L25").Value = -65.497076023391
# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = -57.894736842105
$ws.Range("G26").Value = 67
$ws.Range("H26").Value = -17.910447761194
$ws.Range("I26").Value = 172
$ws.Range("J26").Value = 201
$ws.Range("K26").Value = -14.427860696517
$ws.Range("L26").Value = 22.857142857142
$ws.Range("M26").Value = 43.333333333333
# Row 27
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("L27").Value = -12.5
# Row 28
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = 6.666666666666
$ws.Range("L28").Value = -5.882352941176
# Row 29
$ws.Range("J29").Value = 5
$ws.Range("K29").Value = -40
$ws.Range("L29").Value = -75
$ws.Range("M29").Value = -57.142857142857
# Row 30
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = -25
$ws.Range("L30").Value = -72.727272727272
$ws.Range("M30").Value = -50
